$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet6")
$ws.Range("A1").Value = "test"
